$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared formula in B12 (propagates to C12:AH12 and AG12 joins the shared group)
$ws.Range("B12:AH12").Formula = '=ROUND(AVERAGE(B2:B11),3) &"±"& ROUND(_xlfn.STDEV.P(B2:B11),3)'

# Move the selection (also scrolls the view towards the right-hand columns)
$ws.Range("AG11").Select()

# Adjust column widths (J, AA, AG)
$ws.Columns("J").ColumnWidth = 8.3
$ws.Columns("AA").ColumnWidth = 8.3
$ws.Columns("AG").ColumnWidth = 8
